$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "2" variants of each source (CHP2, solar_th2, pvt2, heat_pump2,
# gas_boiler2) are being dropped from the table, leaving only the "1"
# variants. Those live in columns D, F, H, J and L. Delete them
# right-to-left so the letters of columns not yet processed stay valid.
$ws.Columns("L").Delete()
$ws.Columns("J").Delete()
$ws.Columns("H").Delete()
$ws.Columns("F").Delete()
$ws.Columns("D").Delete()
